$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 475.58334
$ws.Range("J17").Value = 475.58334
$ws.Range("L17").Value = 1426.75002
$ws.Range("N17").Value = -1762.75002

$ws.Range("H112").Value = 2022.1786
$ws.Range("J112").Value = 2325.7827
$ws.Range("L112").Value = 6977.348100000001
$ws.Range("N112").Value = -9193.348100000001

$ws.Range("H137").Value = 1513.3077
$ws.Range("I137").Value = 1470.3636
$ws.Range("J137").Value = 1749.5
$ws.Range("K137").Value = 4411.0908
$ws.Range("L137").Value = 5248.5
$ws.Range("M137").Value = -1861.0908
$ws.Range("N137").Value = -10348.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16674.066
$ws.Range("I32").Value = 19303.244
$ws.Range("J32").Value = 4514.125
$ws.Range("K32").Value = 19303.244
$ws.Range("L32").Value = 4514.125
$ws.Range("M32").Value = -19016.244
$ws.Range("N32").Value = -5088.125

$ws.Range("H61").Value = 2038.138
$ws.Range("I61").Value = 1955.5217
$ws.Range("J61").Value = 2354.8333
$ws.Range("K61").Value = 1955.5217
$ws.Range("L61").Value = 2354.8333
$ws.Range("M61").Value = -1743.5217
$ws.Range("N61").Value = -2778.8333

$ws.Range("H74").Value = 1055.8572
$ws.Range("I74").Value = 650.5263
$ws.Range("K74").Value = 650.5263
$ws.Range("M74").Value = 223.4737

$ws.Range("H77").Value = 1055.8572
$ws.Range("I77").Value = 650.5263
$ws.Range("K77").Value = 3252.6315
$ws.Range("M77").Value = 1115.3685

$ws.Range("H102").Value = 6892.3335
$ws.Range("I102").Value = 7127.5
$ws.Range("J102").Value = 5011
$ws.Range("K102").Value = 7127.5
$ws.Range("L102").Value = 5011
$ws.Range("M102").Value = -5505.5
$ws.Range("N102").Value = -8255

$ws.Range("H132").Value = 7223.9614
$ws.Range("I132").Value = 10198.714
$ws.Range("J132").Value = 3753.4167
$ws.Range("K132").Value = 30596.142
$ws.Range("L132").Value = 11260.2501
$ws.Range("M132").Value = -28066.142
$ws.Range("N132").Value = -16320.2501

$ws.Range("H136").Value = 2038.138
$ws.Range("I136").Value = 1955.5217
$ws.Range("J136").Value = 2354.8333
$ws.Range("K136").Value = 5866.5651
$ws.Range("L136").Value = 7064.499899999999
$ws.Range("M136").Value = -3316.5651
$ws.Range("N136").Value = -12164.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3233.3333
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549

$ws.Range("H103").Value = 49657
$ws.Range("J103").Value = 49657
$ws.Range("L103").Value = 49657
$ws.Range("N103").Value = -52001

$ws.Range("H105").Value = 3306
$ws.Range("I105").Value = 3317.4546
$ws.Range("J105").Value = 3288
$ws.Range("K105").Value = 3317.4546
$ws.Range("L105").Value = 3288
$ws.Range("M105").Value = -1570.4546
$ws.Range("N105").Value = -6782

$ws.Range("H134").Value = 1833.8
$ws.Range("I134").Value = 1559.3871
$ws.Range("K134").Value = 4678.1613
$ws.Range("M134").Value = -2143.1613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2463.6562
$ws.Range("I31").Value = 1746.6818
$ws.Range("J31").Value = 4041
$ws.Range("K31").Value = 1746.6818
$ws.Range("L31").Value = 4041
$ws.Range("M31").Value = -1451.6818
$ws.Range("N31").Value = -4631

$ws.Range("H34").Value = 2463.6562
$ws.Range("I34").Value = 1746.6818
$ws.Range("J34").Value = 4041
$ws.Range("K34").Value = 1746.6818
$ws.Range("L34").Value = 4041
$ws.Range("M34").Value = -1544.6818
$ws.Range("N34").Value = -4445

$ws.Range("H41").Value = 2950
$ws.Range("I41").Value = 2950
$ws.Range("K41").Value = 2950
$ws.Range("M41").Value = -2522

$ws.Range("H51").Value = 10590.75
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4264

$ws.Range("H58").Value = 1373937
$ws.Range("I58").Value = 1765609.6
$ws.Range("J58").Value = 3083
$ws.Range("K58").Value = 1765609.6
$ws.Range("L58").Value = 3083
$ws.Range("M58").Value = -1765406.6
$ws.Range("N58").Value = -3489

$ws.Range("H60").Value = 12730
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 13003
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 13003
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -14025

$ws.Range("H61").Value = 10590.75
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4652

$ws.Range("H132").Value = 542741.6
$ws.Range("I132").Value = 615751.3
$ws.Range("J132").Value = 7337.3335
$ws.Range("K132").Value = 1847253.9
$ws.Range("L132").Value = 22012.0005
$ws.Range("M132").Value = -1844723.9
$ws.Range("N132").Value = -27072.0005

$ws.Range("H134").Value = 1582.9412
$ws.Range("I134").Value = 1232.5676
$ws.Range("J134").Value = 2508.9285
$ws.Range("K134").Value = 3697.7028
$ws.Range("L134").Value = 7526.7855
$ws.Range("M134").Value = -1162.7028
$ws.Range("N134").Value = -12596.7855

$ws.Range("H136").Value = 1373937
$ws.Range("I136").Value = 1765609.6
$ws.Range("J136").Value = 3083
$ws.Range("K136").Value = 5296828.800000001
$ws.Range("L136").Value = 9249
$ws.Range("M136").Value = -5294278.800000001
$ws.Range("N136").Value = -14349

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 11666.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2891.7856
$ws.Range("I132").Value = 2754.6
$ws.Range("J132").Value = 2968
$ws.Range("K132").Value = 8263.799999999999
$ws.Range("L132").Value = 8904
$ws.Range("M132").Value = -5733.799999999999
$ws.Range("N132").Value = -13964

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H49").Value = 50000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H132").Value = 4164.032
$ws.Range("I132").Value = 4232.852
$ws.Range("J132").Value = 3699.5
$ws.Range("K132").Value = 12698.556
$ws.Range("L132").Value = 11098.5
$ws.Range("M132").Value = -10168.556
$ws.Range("N132").Value = -16158.5

$ws.Range("H136").Value = 4313430
$ws.Range("I136").Value = 9617598
$ws.Range("K136").Value = 28852794
$ws.Range("M136").Value = -28850244

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 251749.25

$ws.Range("H84").Value = 251749.25

$ws.Range("H132").Value = 2567.1875
$ws.Range("I132").Value = 2047.3125
$ws.Range("J132").Value = 3087.0625
$ws.Range("K132").Value = 6141.9375
$ws.Range("L132").Value = 9261.1875
$ws.Range("M132").Value = -3611.9375
$ws.Range("N132").Value = -14321.1875

$ws.Range("H136").Value = 2366.7
$ws.Range("I136").Value = 2234.875
$ws.Range("J136").Value = 2454.5833
$ws.Range("K136").Value = 6704.625
$ws.Range("L136").Value = 7363.749899999999
$ws.Range("M136").Value = -4154.625
$ws.Range("N136").Value = -12463.7499
